$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'41.026.47"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -1.48%  "
$ws.Range("D3").Value = "'2.421.83"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -2.10%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "'316.10"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.54%  "
$ws.Range("D6").Value = "'88.75"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -4.13%  "
$ws.Range("D7").Value = "'0.539"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -2.67%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").Value = "'0.495"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -4.62%  "
$ws.Range("D10").Value = "'31.89"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -3.47%  "
$ws.Range("D11").Value = "'0.0829"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -4.47%  "
$ws.Range("E12").Value = "  -2.55%  "
$ws.Range("D13").Value = "'2.798.81"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.94%  "
$ws.Range("D14").Value = "'6.68"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -3.23%  "
$ws.Range("D15").Value = "'15.56"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.18%  "
$ws.Range("D16").Value = "'2.444.45"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.31%  "
$ws.Range("D17").Value = "'0.767"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -2.83%  "
$ws.Range("D18").Value = "'40.961.26"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.54%  "
$ws.Range("D19").Value = "'0.0₃0919"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -3.75%  "
$ws.Range("D20").Value = "'6.23"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -4.07%  "
$ws.Range("D21").Value = "'71.76"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.52%  "
$ws.Range("D22").Value = "'11.02"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -3.18%  "
$ws.Range("D23").Value = "'234.79"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -2.48%  "
$ws.Range("D24").Value = "'2.68"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -2.28%  "
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("E26").Value = "  -3.11%  "
$ws.Range("D27").Value = "'24.01"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.92%  "
$ws.Range("E28").Value = "  -3.55%  "
$ws.Range("D29").Value = "'9.52"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -4.14%  "
$ws.Range("D30").Value = "'34.68"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -4.92%  "
$ws.Range("D31").Value = "'155.79"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.18%  "
$ws.Range("E32").Value = "  +0.05%  "
$ws.Range("D33").Value = "'5.24"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -5.66%  "
$ws.Range("D34").Value = "'2.52"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -2.33%  "
$ws.Range("D35").Value = "'0.0742"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -3.96%  "
$ws.Range("E36").Value = "  +0.80%  "
$ws.Range("E37").Value = "  -4.42%  "
$ws.Range("D38").Value = "'0.114"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.05%  "
$ws.Range("E39").Value = "  -3.95%  "
$ws.Range("D40").Value = "'0.0999"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.04%  "
$ws.Range("D41").Value = "'3.85"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -3.26%  "
$ws.Range("E42").Value = "  -7.61%  "
$ws.Range("D43").Value = "'1.983.49"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.17%  "
$ws.Range("D44").Value = "'18.51"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -3.06%  "
$ws.Range("E45").Value = "  -3.98%  "
$ws.Range("E46").Value = "  -4.78%  "
$ws.Range("D47").Value = "'9.42"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.54%  "
$ws.Range("D48").Value = "'2.660.94"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.83%  "
$ws.Range("D49").Value = "'94.49"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -3.12%  "
$ws.Range("D50").Value = "'72.88"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.90%  "
$ws.Range("D51").Value = "'51.69"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.75%  "
